# B6-PowerPoint.pptx edit
#
# 1) Three tables (on slides 14, 15, 16) get their table style switched
#    from the deck's one custom style ({8561FC4F-442F-4506-A984-9B81CC9B8A49})
#    to the built-in style {E33C04D8-8E35-40AB-915E-1644572D746E}.
#
# 2) The presentation's theme color scheme (backing every slide through the
#    single slide master) is swapped from the "Integral / Red Violet"
#    palette to the "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newStyle = "{E33C04D8-8E35-40AB-915E-1644572D746E}"
$tableSlideIndices = @(14, 15, 16)

foreach ($slideIdx in $tableSlideIndices) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyle)
        }
    }
}

# --- 2. Theme colour scheme -------------------------------------------
# Office palette (dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink), in the same
# order the OOXML <a:clrScheme> children are declared.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$scheme = $p.Slides.Item(1).ThemeColorScheme

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $hex = $officeColors[$i]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgbVal = $r + ($g * 256) + ($b * 65536)
    $scheme.Item($i + 1).RGB = $rgbVal
}
